# Generate Report for Handoff
# Replaces the two e2e test-data file entries (3c88a58f... and 7a7512f2...)
# with fresh ones (27da9dff... and ffffbf64af09...), refreshes the handoff
# status/timestamps, and clears out the stale "Latest Target File" /
# "Latest Handback File" columns on the language sheets.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "3c88a58f-c0ea-4618-8baa-10218af9745c"
$oldGuid2 = "7a7512f2-971c-4baa-9ed6-29d3fcf9ff76"
$newGuid1 = "27da9dff-9113-4d1a-932e-f9925b5fc300"
$newGuid2 = "ffffbf64af09-f01b-42ce-b504-63844fe96fb9"

$newHash = "2788377ae5a81d613a81cd3a36874593a762cf68"

$addrOld1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ad2fa0d2af63346ec30fccc7a8dc5db6cd163cc8/e2e/$oldGuid1.md"
$addrOld2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ad2fa0d2af63346ec30fccc7a8dc5db6cd163cc8/e2e/$oldGuid2.md"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid1.md"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-08 05:21:47"

$wsOverview.Range("A3").Value = "$newGuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newGuid2.md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-08 05:21:47"

# NOTE: Range.Hyperlinks.Delete() removes every hyperlink on the sheet
# (not just the target range), so only call it once per sheet and then
# rebuild every hyperlink that should survive.
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $addrOld1, "", "", "e2e\$newGuid1.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $addrOld2, "", "", "e2e\$newGuid2.md") | Out-Null

$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newGuid1.md"
$wsZh.Range("G2").Value = "$newGuid1.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-09-08 05:21:41"
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

$wsZh.Range("A3").Value = "$newGuid2.md"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = "$newGuid1.$newHash.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-09-08 05:21:41"
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"

# Drop every hyperlink on the sheet (A2, I2, A3, I3) then rebuild only
# A2/A3 - the "Latest Target File" hyperlinks (I2/I3) go away entirely.
$wsZh.Range("A1").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $addrOld1, "", "", "$newGuid1.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $addrOld2, "", "", "$newGuid2.md") | Out-Null

$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("I3").Style = "Normal"

$wsZh.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newGuid1.md"
$wsDe.Range("G2").Value = "$newGuid1.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-09-08 05:21:47"
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDe.Range("A3").Value = "$newGuid2.md"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = "$newGuid1.$newHash.de-de.xlf"
$wsDe.Range("H3").Value = "2016-09-08 05:21:47"
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"

$wsDe.Range("A1").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $addrOld1, "", "", "$newGuid1.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $addrOld2, "", "", "$newGuid2.md") | Out-Null

$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("I3").Style = "Normal"

$wsDe.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426
